$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R (2020 data) mirrors the formatting already used by column Q
# (one year to the left) for the same rows, so copy formats across first.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)  # xlPasteFormats

# Now populate the 2020 values for each indicator row.
$ws.Range("R4").Value = 2020
$ws.Range("R5").Value = 5
$ws.Range("R6").Value = 3.5
$ws.Range("R7").Value = 1.8
$ws.Range("R8").Value = 24.4
$ws.Range("R9").Value = 7.2
$ws.Range("R10").Value = 2.9
$ws.Range("R11").Value = 7.4
$ws.Range("R12").Value = 4
$ws.Range("R13").Value = 3.2
$ws.Range("R14").Value = 3.5

# Match the selection left behind in the saved workbook.
$ws.Range("R4:R14").Select()
